$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17
$ws.Range("A17").Value = "2026-02-15 01:31:50"
$ws.Range("C17").Value = "JULIETTE SOB KEMDJOU"
$ws.Range("D17").Value = 4595

# Row 63
$ws.Range("A63").Value = "2026-02-15 01:13:14"
$ws.Range("D63").Value = 425877

# Row 72
$ws.Range("A72").Value = "2026-02-14 20:16:00"
$ws.Range("D72").Value = 223436

# Row 73
$ws.Range("A73").Value = "2026-02-15 01:38:51"
$ws.Range("C73").Value = "BLANDINE PEYEMBOUO"
$ws.Range("D73").Value = 50822

# Row 74
$ws.Range("A74").Value = "2026-02-15 01:30:06"
$ws.Range("C74").Value = "ETS CAMPUS III ETS MOBILE FINANCIAL SERVICES MFS"
$ws.Range("D74").Value = 104

# Row 76
$ws.Range("A76").Value = "2026-02-15 01:40:35"
$ws.Range("D76").Value = 1098

# Row 78
$ws.Range("A78").Value = "2026-02-14 19:29:43"
$ws.Range("C78").Value = "LAZARRE BIKEK"
$ws.Range("D78").Value = 104496

# Row 85
$ws.Range("A85").Value = "2026-02-15 01:09:36"
$ws.Range("D85").Value = 391865

# Row 87
$ws.Range("A87").Value = "2026-02-15 01:47:52"
$ws.Range("D87").Value = 64097

# Row 89
$ws.Range("A89").Value = "2026-02-15 01:40:24"
$ws.Range("D89").Value = 3842

# Row 94
$ws.Range("A94").Value = "2026-02-15 01:32:29"

$wb.Save()
